$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the heading labels for each solution approach to reflect the
# "one-liner" presentation wording used in the refreshed workbook. B8,
# B13 and B23 were originally quote-prefixed text entries, so a leading
# apostrophe is used to preserve that (and each cell's style/number
# format); B18 was plain text (no quote prefix), so it is set directly.
$ws.Range("B8").Value = "'1) One-liner based on recursion"
$ws.Range("B13").Value = "'2) One-liner based on tail-call recursion"
$ws.Range("B18").Value = "3) One-liner based on array formulas"
$ws.Range("B23").Value = "'4) Use of spreadsheets capabilities"

# Add a closing "(end)" marker below the last data row, entered the same
# quote-prefixed way as the other placeholder cells in column B.
$ws.Range("B46").Value = "'(end)"

# Move the active selection to C3 (previously F3).
$ws.Range("C3").Select()
